$d = $word.ActiveDocument

# Locate the paragraph that carries the "_GoBack" bookmark: in this
# document it is a centered, otherwise-empty paragraph, immediately
# followed by another empty (non-centered) paragraph.
$goBackIndex = 0
for ($i = 1; $i -lt $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    $next = $d.Paragraphs.Item($i + 1)
    if ($pp.Alignment -eq 1 -and $pp.Range.Text.Trim() -eq "" -and $next.Range.Text.Trim() -eq "") {
        $goBackIndex = $i
        break
    }
}
if ($goBackIndex -eq 0) {
    # Fallback: known position in this document.
    $goBackIndex = 5
}

$bookmarkPara = $d.Paragraphs.Item($goBackIndex)

# Step 1: strip the "_GoBack" bookmark out of the centered paragraph,
# leaving it as a plain empty centered paragraph (single-paragraph range
# replaced by single-paragraph XML => clean in-place substitution).
$trimmedXml = @"
<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr></w:pPr></w:p>
"@
$bookmarkPara.Range.InsertXML($trimmedXml) | Out-Null

# Step 2: replace the following empty paragraph with the three new
# paragraphs: a blank spacer, the new "Format instances thread safety"
# list item, and a blank ListParagraph-styled paragraph that now hosts
# the relocated "_GoBack" bookmark.
$newContentXml = @"
<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Format instances thread safety</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">. Remove public static final </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>SimpleDateFormat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> and other Formats everywhere, and replace them with static method that create the instances on demand. We need to do this because Format instances are not thread safe. See </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>io.novaordis.events.api.event</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>.DateProperty</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>getDefaultDateFormat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>) as example.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
"@

$d.Paragraphs.Item($goBackIndex + 1).Range.InsertXML($newContentXml) | Out-Null

Write-Host "Done. Paragraph count now:" $d.Paragraphs.Count
